$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.351.91'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.623.39'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.71%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.11'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.78'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.09%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.849.26'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.626.41'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.04'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.352.59'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.35%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.51'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.75%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '202.96'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.27'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.31'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.05'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.87'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.32'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.46%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.22'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.58'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +9.02%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.19'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.92'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.60%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.157.62'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.43%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.61%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.38'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +4.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.785'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.760.52'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.54%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.64'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +9.99%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.86'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0509'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.81%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.44%  '
